# Rimesso check a inizio_lavorazione anziché fine_lavorazione sulla release date
# per le ricerche locali.
#
# The scheduling rows (A13:R30, 18 data rows) get re-sorted: the job that used
# to land in the last slot (old row 30) now sorts into the first slot (row 13),
# and every other job shifts down by one row (old row N -> new row N+1) for
# N = 13..29. Column values (and their types - numbers vs strings) move as a
# block per row, so read the whole block once and rewrite it rotated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 13
$lastRow = 30
$nRows = $lastRow - $firstRow + 1   # 18
$nCols = 18                          # columns A..R

$srcRange = $ws.Range("A$firstRow`:R$lastRow")
$old = $srcRange.Value2              # 1-indexed [row, col] COM array

# Build the rotated block. Target array for a Value2 assignment must be
# 0-indexed (row 0 = first row of the target range, col 0 = column A).
$new = New-Object 'object[,]' $nRows,$nCols

# New first row (row 13) takes what used to be the last row (row 30).
for ($c = 1; $c -le $nCols; $c++) {
    $new[0, $c - 1] = $old[$nRows, $c]
}

# Every other row shifts down by one: new row r (2..18) = old row r-1.
for ($r = 2; $r -le $nRows; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        $new[$r - 1, $c - 1] = $old[$r - 1, $c]
    }
}

$ws.Range("A$firstRow`:R$lastRow").Value2 = $new
